# Edit script for doc/test/examples/agents.xlsx
# - Replace numeric "Tipo interno" values in column E with descriptive text
#   (Person / Sensor / Entity)
# - Underline-style (no hyperlink) the location value in B4
# - Move the active selection from E8 to E7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric type codes in column E with readable text labels.
$ws.Range("E2").Value = "Person"
$ws.Range("E3").Value = "Person"
$ws.Range("E4").Value = "Person"
$ws.Range("E5").Value = "Sensor"
$ws.Range("E6").Value = "Entity"
$ws.Range("E7").Value = "Sensor"

# Give B4 an underline style (its own font, not a hyperlink color).
$ws.Range("B4").Font.Underline = $true

# Update the active cell/selection to E7 (instead of E8).
$ws.Range("E7").Select()
